# "Clean up form formatting" - unmute_contact.xlsx
#
# The `inputs` group's `source` / `source_id` questions were defined as
# hidden text fields using the (non-standard) XLSForm type "hidden".
# Clean this up to use the standard "text" type together with the
# "hidden" appearance flag instead. The `contact` question's type
# ("string") is likewise normalised to the standard "text" type.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: source
$ws.Range("A4").Value = "text"
$ws.Range("F4").Value = "hidden"

# Row 5: source_id
$ws.Range("A5").Value = "text"
$ws.Range("F5").Value = "hidden"

# Row 7: contact/_id
$ws.Range("A7").Value = "text"

# Leave the cursor where editing finished
$null = $ws.Range("C16").Select()
